$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$data = @(
    @(9,9), @(8,8), @(9,9), @(8,8), @(8,8), @(9,9), @(8,8), @(8,8), @(10,10), @(8,8),
    @(7,7), @(8,8), @(7,7), @(6,7), @(7,7), @(8,8), @(8,8), @(8,8), @(8,8), @(7,7),
    @(8,8), @(8,8), @(6,7), @(8,8), @(8,8), @(7,8), @(8,8), @(8,8), @(8,8), @(7,7),
    @(6,6), @(7,8), @(8,8), @(7,8), @(8,9), @(7,8), @(7,7), @(8,8), @(6,7), @(7,8),
    @(8,8), @(9,9), @(6,7), @(8,8), @(6,6), @(8,9), @(9,9), @(6,6), @(7,8), @(8,8),
    @(7,8), @(5,5), @(8,8), @(8,8), @(9,9), @(7,7), @(8,8), @(9,9), @(8,8), @(8,8),
    @(7,8), @(9,9), @(6,7), @(8,8), @(7,7), @(7,7), @(6,6), @(9,9), @(5,5), @(6,6),
    @(4,4), @(8,9), @(6,7), @(5,5)
)

for ($idx = 0; $idx -lt $data.Length; $idx++) {
    $row = $idx + 2
    $pair = $data[$idx]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
